# Actualización CU y Tajadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows of data (rows 4-6) ---
# Use Value2 with the raw date serial numbers so no ad-hoc date style gets
# auto-applied (we set the real NumberFormat explicitly further below).
$ws.Range("A4").Value2 = 44825
$ws.Range("B4").Value2 = 44825
$ws.Range("C4").Value = "Se genero error al momento del despleigue an el ajuste de calculo de gastos por pruebas automatizadas no contempladas"
$ws.Range("D4").Value = "Se realizo ajuste modificando las pruebas automatizadas para que contemplaran el nuevo campo de gastos"
$ws.Range("E4").Value = "Tito Maturanda - Luis Sabroso"

$ws.Range("A5").Value2 = 44827
$ws.Range("B5").Value2 = 44827
$ws.Range("C5").Value = "Se genera error por que la entidad InformacionAportes no tiene InformacionAportesId embebido al momento de adicionar los ajustes para analisis de codigo SonarQube y Jacoco."
$ws.Range("D5").Value = "Se realiza el ajuste para embeber esta variable para validacion de codigo"
$ws.Range("E5").Value = "Tito Maturanda"

$ws.Range("A6").Value2 = 44828
$ws.Range("B6").Value2 = 44828
$ws.Range("C6").Value = "Se debe generar modificacion del la clase commonDtos por common-dto por validacion de codigo con SpnarQube y Jacoco"
$ws.Range("D6").Value = "Se modifico nombre de la clase para no generar warning en la validacion de codigo."
$ws.Range("E6").Value = "Tito Maturanda"

# --- Date format for the date columns (A2:B6) first, so the numeric cells
#     pick up the workbook's existing built-in date numFmt (14) before any
#     border/alignment gets layered on top. ---
$ws.Range("A2:B6").NumberFormat = "mm-dd-yy"

# --- Borders + wrap text for all data cells (A2:E6), applied last so every
#     cell in a column converges on a single shared style. ---
$dataRange = $ws.Range("A2:E6")
$dataRange.WrapText = $true
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# --- Row heights ---
$ws.Range("2:2").RowHeight = 72.5
$ws.Range("3:3").RowHeight = 72.5
$ws.Range("4:4").RowHeight = 43.5
$ws.Range("5:5").RowHeight = 72.5
$ws.Range("6:6").RowHeight = 43.5

# --- View / selection ---
$ws.Range("C6").Select()
